# Adds a new "2022" data column (M) to the small/medium enterprises table,
# mirroring the existing 2013-2021 columns (D-L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the 2022 column
$ws.Range("M4").Value = 2022
$ws.Range("M5").Value = 2.2
$ws.Range("M6").Value = 1.2

# Match formatting of the neighboring cells:
# - M4 (year header) picks up K4's style (top+bottom border)
$ws.Range("K4").Copy()
$ws.Range("M4").PasteSpecial(-4122) # xlPasteFormats

# - M5 / M6 (data rows) pick up L5 / L6's style (bottom border only)
$ws.Range("L5").Copy()
$ws.Range("M5").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("L6").Copy()
$ws.Range("M6").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

# Move/expand the active sheet's selection to M10, as recorded in the workbook
$ws.Range("M10").Select()
